# feat: add 2022-Q4 data
#
# Inserts a new worksheet "2022-Q4" right after "总计", pushing the
# existing quarterly sheets ("2021-Q3", "2021-Q2", "2021-Q1") one
# position later (their own data is untouched). The "总计" (summary)
# sheet gets a new row for the 2022-Q4 entry and its existing rows are
# relabeled to match the new quarter ordering.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)     # 总计
$wsOldQ3 = $wb.Worksheets.Item(2)     # currently "2021-Q3" - used as a style/structure template

# --- Create the new "2022-Q4" sheet right after "总计" ---------------------
# Duplicating an existing data sheet keeps the same layout/formatting
# (header style, column widths, borders, etc.) as the other quarter sheets.
$wsOldQ3.Copy($null, $wsTotal)

$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"

# Header: "基金金额" -> "基金规模"
$wsQ4.Range("D1").Value = "基金规模"

# Force the text-like columns to stay text (so things like fund codes keep
# leading zeros and numeric-looking values aren't silently converted).
$wsQ4.Range("B2:G2").NumberFormat = "@"

$wsQ4.Range("B2").Value = "233009"
$wsQ4.Range("C2").Value = "大摩多因子精选策略混合"
$wsQ4.Range("D2").Value = "6.42"
$wsQ4.Range("E2").Value = "91.11"
$wsQ4.Range("F2").Value = "0.96"
$wsQ4.Range("G2").Value = "0.0616"
$wsQ4.Range("H2").Value = 6

# --- Update the "总计" summary sheet ---------------------------------------
# Existing row labels shift down one quarter:
#   2021-Q3 -> 2022-Q4 (row 2)
#   2021-Q2 -> 2021-Q3 (row 3)
# and a brand-new row is inserted for 2021-Q2 (row 4), pushing the old
# 2021-Q1 row down to row 5.
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("B3").Value = "2021-Q3"

$wsTotal.Rows.Item(4).Insert()

$wsTotal.Range("A3").Copy($wsTotal.Range("A4"))
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2021-Q2"
$wsTotal.Range("C4").Value = 1
$wsTotal.Range("D4").Value = 0.06

$wsTotal.Range("A5").Value = 3

# Keep "总计" as the active sheet/tab, matching the original view state.
$wsTotal.Activate()
